$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Export Demand" bound table is being extended with two extra year rows
# (2035 and 2040, using the same low bound as 2030) while the existing
# 2045/2050 rows are pushed down two rows to make room for them.

# --- Prepare formatting for the two newly-populated rows (8 and 9) ---
# These rows were blank spacer rows before; give their data columns the
# same look as the existing data rows (6/7) before filling in values.
$ws.Range("D8:F8").ClearFormats()
$ws.Range("D9:F9").ClearFormats()
$ws.Range("G6:M6").Copy($ws.Range("G8:M8"))
$ws.Range("G6:M6").Copy($ws.Range("G9:M9"))

# --- Row 6: was Year 2045 / DKW 50, becomes Year 2035 / DKW 10 ---
$ws.Range("D6").Value = 2035
$ws.Range("E6").Value = "LO"
$ws.Range("F6").Value = "ACT_BND"
$ws.Range("H6").Value = 10
$ws.Range("J6").Value = 0
$ws.Range("M6").Value = "EXPH2*"

# --- Row 7: was Year 2050 / DKW 50, becomes Year 2040 / DKW 10 ---
$ws.Range("D7").Value = 2040
$ws.Range("E7").Value = "LO"
$ws.Range("F7").Value = "ACT_BND"
$ws.Range("H7").Value = 10
$ws.Range("J7").Value = 0
$ws.Range("M7").Value = "EXPH2*"

# --- Row 8 (new data row): Year 2045 / DKW 50 ---
$ws.Range("D8").Value = 2045
$ws.Range("E8").Value = "LO"
$ws.Range("F8").Value = "ACT_BND"
$ws.Range("H8").Value = 50
$ws.Range("J8").Value = 0
$ws.Range("M8").Value = "EXPH2*"

# --- Row 9 (new data row): Year 2050 / DKW 50 ---
$ws.Range("D9").Value = 2050
$ws.Range("E9").Value = "LO"
$ws.Range("F9").Value = "ACT_BND"
$ws.Range("H9").Value = 50
$ws.Range("J9").Value = 0
$ws.Range("M9").Value = "EXPH2*"

# Restore the selection that was active when the sheet was last saved.
$ws.Range("O6").Select()
